$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07967938146308559
$ws.Range("C2").Value = 0.8019381115720245
$ws.Range("D2").Value = 1.602517945319166
$ws.Range("E2").Value = 1.265905978072292
$ws.Range("F2").Value = 1.275721982108824
$ws.Range("B3").Value = 0.1996680100171121
$ws.Range("C3").Value = 1.350499864034121
$ws.Range("D3").Value = 3.617881843527005
$ws.Range("E3").Value = 1.902073038431228
$ws.Range("F3").Value = 1.910386039492769
$ws.Range("B4").Value = 0.6016847796089962
$ws.Range("C4").Value = 1.493953804328142
$ws.Range("D4").Value = 4.491473340817588
$ws.Range("E4").Value = 2.119309637787171
$ws.Range("F4").Value = 2.052735548343508
$ws.Range("B5").Value = 0.3601685393389483
$ws.Range("C5").Value = 1.499852406935977
$ws.Range("D5").Value = 4.026777908913376
$ws.Range("E5").Value = 2.006683310568306
$ws.Range("F5").Value = 1.994553887448934
$ws.Range("G5").Value = 49
$ws.Range("B6").Value = 0.554109530473832
$ws.Range("C6").Value = 1.40221327316157
$ws.Range("D6").Value = 3.723570466260433
$ws.Range("E6").Value = 1.929655530466625
$ws.Range("F6").Value = 1.867946810805509
$ws.Range("G6").Value = 48
$ws.Range("B7").Value = 0.3968503597739336
$ws.Range("C7").Value = 1.194979168351206
$ws.Range("D7").Value = 2.489467468001566
$ws.Range("E7").Value = 1.5778046355622
$ws.Range("F7").Value = 1.547043982108814
$ws.Range("G7").Value = 39
$ws.Range("B8").Value = 0.5649144946704114
$ws.Range("C8").Value = 1.22464266767186
$ws.Range("D8").Value = 2.755048971642472
$ws.Range("E8").Value = 1.659834019305085
$ws.Range("F8").Value = 1.581694116082434
$ws.Range("G8").Value = 38
$ws.Range("B9").Value = 0.884415467835616
$ws.Range("C9").Value = 1.104993751378739
$ws.Range("D9").Value = 2.323198242346433
$ws.Range("E9").Value = 1.524204134079957
$ws.Range("F9").Value = 1.27202904791106
$ws.Range("G9").Value = 21
$ws.Range("B10").Value = 0.5881942522924533
$ws.Range("C10").Value = 0.8478029123914703
$ws.Range("D10").Value = 0.9561329667000577
$ws.Range("E10").Value = 0.9778205186536318
$ws.Range("F10").Value = 0.810614526421041
$ws.Range("G10").Value = 14
$ws.Range("B11").Value = 0.7033621104908662
$ws.Range("C11").Value = 0.8247366046038167
$ws.Range("D11").Value = 1.023412291882237
$ws.Range("E11").Value = 1.011638419536465
$ws.Range("F11").Value = 0.812937600163807
